{"js": "// Edit the \"Jefe del \u00c1rea de Embarque\" user-story paragraph:\n//   - \"tener un control sobre\" -> \"almacenar\"\n//   - append \" y que hayan sido validados por la entidad correspondiente\"\n//     right after \"...Informe T\u00e9cnico de Vuelo) \" (before \"para ...\")\n//\n// Both edits are plain text replacements inside runs that keep the\n// surrounding (bold, accent6/984806-colored) character formatting, which\n// Word does automatically when a search hit is replaced in place.\n\nconst body = context.document.body;\n\n// 1) \"tener un control sobre\" -> \"almacenar\"\nconst targetPhrase = body.search(\"tener un control sobre\", { matchCase: true });\ntargetPhrase.load(\"items,text\");\n\n// 2) Insert the new clause right after the closing parenthesis clause,\n//    before \"para que est\u00e9 de acuerdo...\"\nconst closingClause = body.search(\n  \"Plan de Vuelo, Control de Peso y Balance, Informe T\u00e9cnico de Vuelo) \",\n  { matchCase: true }\n);\nclosingClause.load(\"items,text\");\n\nawait context.sync();\n\nif (targetPhrase.items.length === 0) {\n  throw new Error('Could not find \"tener un control sobre\" in the document body.');\n}\nif (closingClause.items.length === 0) {\n  throw new Error('Could not find the \"...Informe T\u00e9cnico de Vuelo) \" clause in the document body.');\n}\n\n// Replace \"tener un control sobre\" with \"almacenar\" \u2014 in-place replace keeps\n// the run's existing bold/orange (accent6) formatting.\ntargetPhrase.items[0].insertText(\"almacenar\", Word.InsertLocation.replace);\n\n// Re-insert the closing-parenthesis clause with the new trailing text added,\n// preserving the same trailing single space before \"para\".\nclosingClause.items[0].insertText(\n  \"Plan de Vuelo, Control de Peso y Balance, Informe T\u00e9cnico de Vuelo) y que hayan sido validados por la entidad correspondiente \",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Edit the \"Jefe del \u00c1rea de Embarque\" user-story paragraph:\n#   - \"tener un control sobre\" -> \"almacenar\"\n#   - append \" y que hayan sido validados por la entidad correspondiente\"\n#     right after \"...Informe T\u00e9cnico de Vuelo) \" (before \"para ...\")\n#\n# Both edits are done with Find/Replace on $d.Content, which rewrites the\n# matched text in place and keeps the surrounding (bold, accent6/984806)\n# character formatting of the runs it touches.\n\n$d = $word.ActiveDocument\n\n# 1) \"tener un control sobre\" -> \"almacenar\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"tener un control sobre\"\n$find1.Replacement.Text = \"almacenar\"\n$find1.Execute(\n    $find1.Text,    # FindText\n    $false,         # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap (wdFindContinue)\n    $false,         # Format\n    $find1.Replacement.Text,  # ReplaceWith\n    2               # Replace (wdReplaceOne)\n) | Out-Null\n\n# 2) Insert the new clause right after the closing parenthesis clause,\n#    before \"para que est\u00e9 de acuerdo...\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Informe T\u00e9cnico de Vuelo) \"\n$find2.Replacement.Text = \"Informe T\u00e9cnico de Vuelo) y que hayan sido validados por la entidad correspondiente \"\n$find2.Execute(\n    $find2.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find2.Replacement.Text,\n    2\n) | Out-Null\n"}
